# Update "想去人数" (want-to-go count) figures that changed between
# crawler runs, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6535
$ws1.Range("F5").Value = 408
$ws1.Range("F12").Value = 162
$ws1.Range("F14").Value = 1113
$ws1.Range("F15").Value = 3233
$ws1.Range("F18").Value = 1887

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 3

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6535
$ws4.Range("F5").Value = 408
$ws4.Range("F8").Value = 3
$ws4.Range("F13").Value = 162
$ws4.Range("F15").Value = 1113
$ws4.Range("F16").Value = 3233
$ws4.Range("F19").Value = 1887
